# Update Config.xlsx with Input/Output file path
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Settings")

# Add new configuration rows (Name / Value) for the git-diff bot input & output paths.
$ws.Range("A6").Value = "InputCSVFilePath"
$ws.Range("B6").Value = "C:\Users\angel\GitHub\uipath-automation-4\ProjectPlagiarismBot\Data\Input\GitHubRepoURLInput.CSV"

$ws.Range("A7").Value = "GitCloneRootFilePath"
$ws.Range("B7").Value = "C:\Users\angel\GitHub\uipath-automation-4\ProjectPlagiarismBot\GitRepoCloning"

$ws.Range("A8").Value = "GitDiffOutputFilePath"
$ws.Range("B8").Value = "C:\Users\angel\GitHub\uipath-automation-4\ProjectPlagiarismBot\Data\Output"

# Make Settings the active sheet/tab again (it was Assets before the edit).
$ws.Activate()
$ws.Range("A1").Select()
